# The table's last column (column 5 - the "p" / significance-test-name
# column, holding values like "test"/"exact") is removed entirely: every
# row loses its last cell.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$tbl.Columns.Item(5).Delete()
